$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.68
$ws.Range("H3").Value = 2.66
$ws.Range("K3").Value = 4.2
$ws.Range("L3").Value = 1.31
$ws.Range("N3").Value = 3.9
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 2.02
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 1.39
$ws.Range("S3").Value = 3.15
$ws.Range("T3").Value = 1.7
$ws.Range("U3").Value = 2.2
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 23
$ws.Range("AF3").Value = 23
$ws.Range("AN3").Value = 27
$ws.Range("AO3").Value = 28
$ws.Range("G4").Value = 2.68
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 4.1
$ws.Range("Q4").Value = 1.74
$ws.Range("S4").Value = 2.86
$ws.Range("T4").Value = 1.64
$ws.Range("W4").Value = 1.59
$ws.Range("F5").Value = 8.4
$ws.Range("H5").Value = 1.43
$ws.Range("I5").Value = 1.44
$ws.Range("K5").Value = 5.4
$ws.Range("L5").Value = 1.33
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 5.1
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 2.42
$ws.Range("Q5").Value = 1.68
$ws.Range("R5").Value = 1.55
$ws.Range("S5").Value = 2.74
$ws.Range("V5").Value = 3.25
$ws.Range("X5").Value = 23
$ws.Range("Y5").Value = 9.4
$ws.Range("AB5").Value = 30
$ws.Range("AC5").Value = 11.5
$ws.Range("AE5").Value = 14
$ws.Range("AL5").Value = 100
$ws.Range("AN5").Value = 120
$ws.Range("AO5").Value = 5.8
$ws.Range("H7").Value = 1.88
$ws.Range("I7").Value = 1.89
$ws.Range("AH7").Value = 17.5
$ws.Range("AI7").Value = 30
$ws.Range("AK7").Value = 50
$ws.Range("AM7").Value = 85
$ws.Range("F9").Value = 2.36
$ws.Range("O9").Value = 1.34
$ws.Range("P9").Value = 1.91
$ws.Range("AB9").Value = 9.8
$ws.Range("O10").Value = 1.24
$ws.Range("F11").Value = 2.76
$ws.Range("G11").Value = 2.78
$ws.Range("H11").Value = 2.84
$ws.Range("I11").Value = 2.86
$ws.Range("L11").Value = 1.42
$ws.Range("P11").Value = 1.99
$ws.Range("V11").Value = 1.53
$ws.Range("W11").Value = 1.56
$ws.Range("Z11").Value = 18
$ws.Range("AK11").Value = 29
$ws.Range("AN11").Value = 24
$ws.Range("Q12").Value = 1.7
$ws.Range("AO12").Value = 10
$ws.Range("G13").Value = 4.9
$ws.Range("R13").Value = 1.65
$ws.Range("T13").Value = 1.61
$ws.Range("H14").Value = 4.1
$ws.Range("N14").Value = 3.05
$ws.Range("O14").Value = 1.43
$ws.Range("R14").Value = 1.26
$ws.Range("S14").Value = 4.2
$ws.Range("U14").Value = 1.9
$ws.Range("V14").Value = 1.28
